$wb = $excel.ActiveWorkbook

# --- Sheet 1 "PI hours": add a new "cfop" column (G) ---------------------
$ws1 = $wb.Worksheets.Item("PI hours")

# Header cell, copy formatting (bold + border) from an existing header cell
$headerSrc = $ws1.Range("B1")
$headerSrc.Copy()
$ws1.Range("G1").PasteSpecial(-4122)  # xlPasteFormats
$ws1.Range("G1").Value = "cfop"

# Per-row cfop application lists
$ws1.Range("G2").Value = "['cfop_GC']"
$ws1.Range("G3").Value = "['cfop_NH']"
$ws1.Range("G4").Value = "['cfop_SHA']"

$excel.CutCopyMode = $false

# --- New sheet "cfop hours" -----------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add($null, $lastSheet)
$ws4.Name = "cfop hours"

# Copy the layout/formatting of the "department hours" sheet (same shape)
$deptSheet = $wb.Worksheets.Item("department hours")
$deptSheet.Range("B1:D1").Copy()
$ws4.Range("B1:D1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$deptSheet.Range("A2:A4").Copy()
$ws4.Range("A2:A4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws4.Range("B1").Value = "cfop"
$ws4.Range("C1").Value = "hours"
$ws4.Range("D1").Value = "percentage"

$ws4.Range("A2").Value = 0
$ws4.Range("B2").Value = "cfop_GC"
$ws4.Range("C2").Value = 43
$ws4.Range("D2").Value = 67.71653543307086

$ws4.Range("A3").Value = 1
$ws4.Range("B3").Value = "cfop_NH"
$ws4.Range("C3").Value = 14.5
$ws4.Range("D3").Value = 22.83464566929134

$ws4.Range("A4").Value = 2
$ws4.Range("B4").Value = "cfop_SHA"
$ws4.Range("C4").Value = 6
$ws4.Range("D4").Value = 9.448818897637794

# Keep original active sheet/tab selected as it was before the edit
$ws1.Activate()
